$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column CS (18-sep) with header + data
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy formatting (bold / border / centered) of the last header cell (CR1)
# onto the new header cell (CS1), then set its text.
$wsSpot.Range("CR1").Copy($wsSpot.Range("CS1"))
$wsSpot.Range("CS1").Value = "18-sep"

# Fill in the 24 hourly values for the new day.
$wsSpot.Range("CS2").Value = 27.6
$wsSpot.Range("CS3").Value = 27.02
$wsSpot.Range("CS4").Value = 27.4
$wsSpot.Range("CS5").Value = 15.61
$wsSpot.Range("CS6").Value = 17.23
$wsSpot.Range("CS7").Value = 25.65
$wsSpot.Range("CS8").Value = 40.7
$wsSpot.Range("CS9").Value = 100.05
$wsSpot.Range("CS10").Value = 102.63
$wsSpot.Range("CS11").Value = 79.09
$wsSpot.Range("CS12").Value = 56.4
$wsSpot.Range("CS13").Value = 5.87
$wsSpot.Range("CS14").Value = 0.2
$wsSpot.Range("CS15").Value = 0
$wsSpot.Range("CS16").Value = 0
$wsSpot.Range("CS17").Value = 0.52
$wsSpot.Range("CS18").Value = 15.31
$wsSpot.Range("CS19").Value = 45.8
$wsSpot.Range("CS20").Value = 96.86
$wsSpot.Range("CS21").Value = 124.7
$wsSpot.Range("CS22").Value = 127.5
$wsSpot.Range("CS23").Value = 110
$wsSpot.Range("CS24").Value = 99.09
$wsSpot.Range("CS25").Value = 88.2

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new day row (2025-09-16 / 31.925)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the new date cell to stay a plain text label (matching the rest
# of column A) instead of being auto-converted into a date serial value.
$wsGaz.Range("A94").NumberFormat = "@"
$wsGaz.Range("A94").Value = "2025-09-16"
$wsGaz.Range("B94").Value = 31.925
# Re-apply the (unstyled) format of the previous row so the new cell
# doesn't keep a "text" number format applied to it.
$wsGaz.Range("A93").Copy()
$wsGaz.Range("A94").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet "CO2": append a new day row (2025-09-16 / 77.29)
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A94").NumberFormat = "@"
$wsCo2.Range("A94").Value = "2025-09-16"
$wsCo2.Range("B94").Value = 77.29000000000001
$wsCo2.Range("A93").Copy()
$wsCo2.Range("A94").PasteSpecial(-4122)
